$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Remove the previous test results. The KMeans evaluation now returns
#    the first section of each cluster instead of pre-computed numbers, so
#    the old cached measurements for the Electron proxy (B/C) and No proxy
#    (E/F) columns are cleared, leaving just the iteration index in A.
# ---------------------------------------------------------------------------
$ws.Range("B10:C19").ClearContents()
$ws.Range("E10:F19").ClearContents()

# The Electron proxy averages (B20/C20) no longer exist at all now that
# their source data is gone.
$ws.Range("B20:C20").ClearContents()

# ---------------------------------------------------------------------------
# 2. Append a third copy of the "cluster" table (same layout as the table
#    that already lives at rows 23-43) further down the sheet, starting at
#    row 46.
# ---------------------------------------------------------------------------

# Blank spacer block (same look as A1:F6 / A23:F29) -------------------------
$ws.Range("A46:F52").HorizontalAlignment = -4108
$ws.Range("A46:F52").WrapText = $true

# Header row (row 54) - mirrors row 8 / row 31 -------------------------------
$ws.Range("A54").Value2 = "Iteration"
$ws.Range("B54").Value2 = "Electron Proxy - BBC"
$ws.Range("E54").Value2 = "No Proxy - BBC"
$ws.Range("H54").Value2 = "System Proxy - BBC"

$ws.Range("B54:C54").HorizontalAlignment = -4108
$ws.Range("B54:C54").WrapText = $true
$ws.Range("D54").WrapText = $true
$ws.Range("E54:F54").HorizontalAlignment = -4108
$ws.Range("H54:I54").HorizontalAlignment = -4108
$ws.Range("H54:I54").WrapText = $true

# Sub-header row (row 55) - mirrors row 9 / row 32 ---------------------------
$ws.Range("B55").Value2 = "First Request"
$ws.Range("C55").Value2 = "Second Request"
$ws.Range("E55").Value2 = "First Request"
$ws.Range("F55").Value2 = "Second Request"
$ws.Range("H55").Value2 = "First Request"
$ws.Range("I55").Value2 = "Second Request"

# Iteration index column (rows 56-65) - mirrors rows 10-19 / 33-42 ----------
$ws.Range("A56").Value2 = 1
$ws.Range("A57").Value2 = 2
$ws.Range("A58").Value2 = 3
$ws.Range("A59").Value2 = 4
$ws.Range("A60").Value2 = 5
$ws.Range("A61").Value2 = 6
$ws.Range("A62").Value2 = 7
$ws.Range("A63").Value2 = 8
$ws.Range("A64").Value2 = 9
$ws.Range("A65").Value2 = 10

# Averages row (row 66) - mirrors row 20 / row 43 ----------------------------
$ws.Range("A66").Value2 = "Avg"
$ws.Range("B66").Formula = "=AVERAGE(B56:B65)"
$ws.Range("C66").Formula = "=AVERAGE(C56:C65)"
$ws.Range("E66").Formula = "=AVERAGE(E56:E65)"
$ws.Range("F66").Formula = "=AVERAGE(F56:F65)"
$ws.Range("H66").Formula = "=AVERAGE(H56:H65)"
$ws.Range("I66").Formula = "=AVERAGE(I56:I65)"

# ---------------------------------------------------------------------------
# 3. Merge the new header cells the same way the other two tables do.
# ---------------------------------------------------------------------------
$ws.Range("A46:F52").Merge()
$ws.Range("B54:C54").Merge()
$ws.Range("E54:F54").Merge()
$ws.Range("H54:I54").Merge()

# ---------------------------------------------------------------------------
# 4. Update the view state to match where the user ended up after the edit.
# ---------------------------------------------------------------------------
$ws.Activate()
try {
    $ws.Range("A43").Select()
    $excel.ActiveWindow.ScrollRow = 43
} catch {
}
$ws.Range("R57").Select()
